# Fix the harvester column: replace "Retrofitted_480" placeholder values
# with the correct harvester initials "S.GISH" (holly added S.GISH to
# harvester in bioSamples).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2:B13").Value = "S.GISH"

# The last row (14) is completely empty - drop it so the sheet's used
# range / dimension shrinks back down to row 13.
$ws.Rows("14").Delete()

# Nudge column B a bit wider to fit the new harvester values, and select
# the column (mirrors the manual edit that produced this change).
$ws.Columns("B").ColumnWidth = 8
$ws.Columns("B").Select() | Out-Null

# Re-apply bold to the rnaPrepMethod header so it collapses onto the same
# (already existing) bold style as the rest of the header row instead of
# keeping its own redundant font/style entry.
$ws.Range("G1").Font.Bold = $true
